$wb = $excel.ActiveWorkbook

# --- Sheet "leaderboard2" (Cobblemons) ---
$ws1 = $wb.Worksheets.Item("leaderboard2")
$ws1.Range("D3").Value = 312
$ws1.Range("D4").Value = 121
$ws1.Range("B13").Value = "Dernière update le 27.03.25 à 23:00"

# --- Sheet "leaderboard3" (Shiny Cobblemons) ---
$ws2 = $wb.Worksheets.Item("leaderboard3")
$ws2.Range("C3").Value = "BKZRackham"
$ws2.Range("D3").Value = 25
$ws2.Range("C4").Value = "ArtyumsM"
$ws2.Range("D4").Value = 20
$ws2.Range("B13").Value = "Dernière update le 27.03.25 à 23:00"
